# 7.1.2.xlsx update:
#  - add a new "2022" year column (H) to every data row that currently
#    carries 2018-2021 figures
#  - insert a new "By sex" / "Men" / "Women" block (3 rows) right after
#    the "Rural/Urban" block and before the "By territory" block
#  - move the selection to I3 (mirrors the author's last-saved cursor)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the 3 new rows for the "By sex" block above the old row 8
#    ("By territory" header), pushing everything below it down by 3.
# ---------------------------------------------------------------------
$ws.Rows("8:10").Insert()

# ---------------------------------------------------------------------
# 2) New row 8: section header "Жынысы боюнча" / "По полу" / "By sex"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Жынысы боюнча"
$ws.Range("B8").Value = "По полу"
$ws.Range("C8").Value = "By sex"

$ws.Range("A8:C8").Font.Name = "Times New Roman"
$ws.Range("A8:C8").Font.Size = 9
$ws.Range("A8:C8").Font.Bold = $true
$ws.Range("A8:C8").Font.Italic = $true
$ws.Range("A8:C8").HorizontalAlignment = -4131
$ws.Range("A8:C8").VerticalAlignment = -4160
$ws.Range("A8:C8").WrapText = $true
$ws.Range("A8:C8").IndentLevel = 1

$ws.Range("D8").Font.Name = "Times New Roman"
$ws.Range("D8").Font.Size = 9
$ws.Range("D8").Font.Bold = $true
$ws.Range("D8").NumberFormat = "0.0"
$ws.Range("D8").HorizontalAlignment = -4152
$ws.Range("D8").VerticalAlignment = -4108

$ws.Range("F8").Font.Name = "Times New Roman"
$ws.Range("F8").Font.Size = 9
$ws.Range("F8").Font.Bold = $true
$ws.Range("F8").NumberFormat = "0.0"
$ws.Range("F8").HorizontalAlignment = -4152
$ws.Range("F8").VerticalAlignment = -4108

$ws.Range("E8,G8,H8").Font.Name = "Times New Roman"
$ws.Range("E8,G8,H8").Font.Size = 9
$ws.Range("E8,G8,H8").Font.Bold = $true
$ws.Range("E8,G8,H8").NumberFormat = "0.0"

# ---------------------------------------------------------------------
# 3) New row 9: "Эркектер" / "Мужчины" / "Men" (male figures)
# ---------------------------------------------------------------------
$ws.Range("A9").Value = "Эркектер"
$ws.Range("B9").Value = "Мужчины"
$ws.Range("C9").Value = "Men"

$ws.Range("A9:C9").Font.Name = "Times New Roman"
$ws.Range("A9:C9").Font.Size = 9
$ws.Range("A9:C9").HorizontalAlignment = -4131
$ws.Range("A9:C9").VerticalAlignment = -4160
$ws.Range("A9:C9").WrapText = $true
$ws.Range("A9:C9").IndentLevel = 1

$ws.Range("D9").Value = "-"
$ws.Range("E9").Value = 20.4
$ws.Range("F9").Value = 20.5
$ws.Range("G9").Value = 19
$ws.Range("H9").Value = 20.145701762391958

$ws.Range("D9,F9").Font.Name = "Times New Roman"
$ws.Range("D9,F9").Font.Size = 9
$ws.Range("D9,F9").NumberFormat = "0.0"
$ws.Range("D9,F9").HorizontalAlignment = -4152
$ws.Range("D9,F9").VerticalAlignment = -4108

$ws.Range("E9,G9,H9").Font.Name = "Times New Roman"
$ws.Range("E9,G9,H9").Font.Size = 9
$ws.Range("E9,G9,H9").NumberFormat = "0.0"

# ---------------------------------------------------------------------
# 4) New row 10: "Аялдар" / "Женщины" / "Women" (female figures)
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "Аялдар"
$ws.Range("B10").Value = "Женщины"
$ws.Range("C10").Value = "Women"

$ws.Range("A10:C10").Font.Name = "Times New Roman"
$ws.Range("A10:C10").Font.Size = 9
$ws.Range("A10:C10").HorizontalAlignment = -4131
$ws.Range("A10:C10").VerticalAlignment = -4160
$ws.Range("A10:C10").WrapText = $true
$ws.Range("A10:C10").IndentLevel = 1

$ws.Range("D10").Value = "-"
$ws.Range("E10").Value = 23.8
$ws.Range("F10").Value = 23.6
$ws.Range("G10").Value = 22
$ws.Range("H10").Value = 22.813769684575334

$ws.Range("D10,F10").Font.Name = "Times New Roman"
$ws.Range("D10,F10").Font.Size = 9
$ws.Range("D10,F10").NumberFormat = "0.0"
$ws.Range("D10,F10").HorizontalAlignment = -4152
$ws.Range("D10,F10").VerticalAlignment = -4108

$ws.Range("E10,G10,H10").Font.Name = "Times New Roman"
$ws.Range("E10,G10,H10").Font.Size = 9
$ws.Range("E10,G10,H10").NumberFormat = "0.0"

# ---------------------------------------------------------------------
# 5) New "2022" column (H) for every row that already has 2018-2021 data
# ---------------------------------------------------------------------
$ws.Range("H3").Value = 2022
$ws.Range("H3").NumberFormat = $ws.Range("G3").NumberFormat
$ws.Range("H3").Font.Name = $ws.Range("G3").Font.Name
$ws.Range("H3").Font.Size = $ws.Range("G3").Font.Size
$ws.Range("H3").Font.Bold = $ws.Range("G3").Font.Bold
$ws.Range("H3").HorizontalAlignment = $ws.Range("G3").HorizontalAlignment
$ws.Range("H3").VerticalAlignment = $ws.Range("G3").VerticalAlignment

function Copy-ColumnFormat($srcAddr, $dstAddr) {
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)
    $dst.NumberFormat = $src.NumberFormat
    $dst.Font.Name = $src.Font.Name
    $dst.Font.Size = $src.Font.Size
    $dst.Font.Bold = $src.Font.Bold
    $dst.Font.Italic = $src.Font.Italic
    $dst.HorizontalAlignment = $src.HorizontalAlignment
    $dst.VerticalAlignment = $src.VerticalAlignment
}

$ws.Range("H4").Value = 21.549331200908018
Copy-ColumnFormat "G4" "H4"

Copy-ColumnFormat "G5" "H5"

$ws.Range("H6").Value = 52.326989300763088
Copy-ColumnFormat "G6" "H6"

$ws.Range("H7").Value = 3.6916333239218613
Copy-ColumnFormat "G7" "H7"

$ws.Range("H12").Value = 6.3300243735913453
Copy-ColumnFormat "G12" "H12"

$ws.Range("H13").Value = 14.956994359947711
Copy-ColumnFormat "G13" "H13"

$ws.Range("H14").Value = 10.216783741088173
Copy-ColumnFormat "G14" "H14"

$ws.Range("H15").Value = 12.292586981974505
Copy-ColumnFormat "G15" "H15"

$ws.Range("H16").Value = 1.1904474908258234
Copy-ColumnFormat "G16" "H16"

$ws.Range("H17").Value = 13.453576190228787
Copy-ColumnFormat "G17" "H17"

$ws.Range("H18").Value = 16.392955140305855
Copy-ColumnFormat "G18" "H18"

$ws.Range("H19").Value = 73.546822306129386
Copy-ColumnFormat "G19" "H19"

$ws.Range("H20").Value = 33.802493306724259
Copy-ColumnFormat "G20" "H20"

$ws.Range("H22").Value = "-"
Copy-ColumnFormat "G22" "H22"

$ws.Range("H23").Value = "-"
Copy-ColumnFormat "G23" "H23"

$ws.Range("H24").Value = "-"
Copy-ColumnFormat "G24" "H24"

$ws.Range("H25").Value = "-"
Copy-ColumnFormat "G25" "H25"

$ws.Range("H26").Value = "-"
Copy-ColumnFormat "G26" "H26"

Copy-ColumnFormat "G27" "H27"

$ws.Range("H28").Value = "-"
Copy-ColumnFormat "G28" "H28"

$ws.Range("H29").Value = "-"
Copy-ColumnFormat "G29" "H29"

$ws.Range("H30").Value = "-"
Copy-ColumnFormat "G30" "H30"

$ws.Range("H31").Value = "-"
Copy-ColumnFormat "G31" "H31"

$ws.Range("H32").Value = "-"
Copy-ColumnFormat "G32" "H32"

# ---------------------------------------------------------------------
# 6) Restore the cursor position the author last left the sheet in
# ---------------------------------------------------------------------
$ws.Range("I3").Select() | Out-Null
